$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Column E ("ok") / Column F (various comments) additions -------------

$ok = "ok"
$h2 = "H-2 is closest available"
$h2_15 = "H-2 is closest available in 15min resolution. However, there is one-minute imbalance price of H-1"
$whatIsThis = "What is this? "
$okIThink = "ok, I think"
$upToH4 = "available up to H-4"
$notFound = "Not found"

# Crossborder flow block (rows 2-6): ok + H-2 is closest available
$ws.Range("E2").Value = $ok
$ws.Range("F2").Value = $h2
$ws.Range("E3").Value = $ok
$ws.Range("F3").Value = $h2
$ws.Range("E4").Value = $ok
$ws.Range("F4").Value = $h2
$ws.Range("E5").Value = $ok
$ws.Range("F5").Value = $h2
$ws.Range("E6").Value = $ok
$ws.Range("F6").Value = $h2

# Generation block (rows 7-9): ok only
$ws.Range("E7").Value = $ok
$ws.Range("E8").Value = $ok
$ws.Range("E9").Value = $ok

# SI and imbalance price block (rows 10-17): ok + detail comment (row 12 differs)
$ws.Range("E10").Value = $ok
$ws.Range("F10").Value = $h2_15
$ws.Range("E11").Value = $ok
$ws.Range("F11").Value = $h2_15
$ws.Range("E12").Value = $ok
$ws.Range("F12").Value = $whatIsThis
$ws.Range("E13").Value = $ok
$ws.Range("F13").Value = $h2_15
$ws.Range("E14").Value = $ok
$ws.Range("F14").Value = $h2_15
$ws.Range("E15").Value = $ok
$ws.Range("F15").Value = $h2_15
$ws.Range("E16").Value = $ok
$ws.Range("F16").Value = $h2_15
$ws.Range("E17").Value = $ok
$ws.Range("F17").Value = $h2_15

# Load block (rows 18-19): only F column comments
$ws.Range("F18").Value = $upToH4
$ws.Range("F19").Value = $notFound

# RES forecast block (rows 20-22): ok only
$ws.Range("E20").Value = $ok
$ws.Range("E21").Value = $ok
$ws.Range("E22").Value = $ok

# ARC merit order row (28): ok, I think
$ws.Range("E28").Value = $okIThink

# --- View state: zoom + selection -----------------------------------------
$ws.Range("F23").Select()
$excel.ActiveWindow.Zoom = 73
